# Apply the commit: "updated sampling ranges for experiment and added climate change factor to hydropower"
# This inserts a new row (row 4) into the "strategy_id-0" sheet for the variable
# climate_change_factor_gnrl_hydropower_availability, shifting the existing
# General-subsector rows (elasticity_gnrl_rate_occupancy_to_gdppc ... population_gnrl_urban)
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Insert a brand new row above current row 4 (pushes rows 4..11 down to 5..12)
$ws.Rows.Item(4).Insert()

# Fill in the new row 4 with the climate change factor variable.
$ws.Cells.Item(4, 1).Value = "General"
$ws.Cells.Item(4, 2).Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 0.5

# Columns J (10) through AS (45) are all set to 1 for this new row.
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(4, $col).Value = 1
}
